$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Total for row 4 (EEID 1234, John Doe) changes from 42.34 to 54.34
$ws.Range("C4").Value = 54.34

# A new employee record ("Jake Doe") is inserted as row 8, pushing the
# existing rows 8 (Joe Doe) and 9 (Irvin Doe) down to rows 9 and 10.
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 5555
$ws.Range("B8").Value = "Jake Doe"
$ws.Range("C8").Value = 21.5
$ws.Range("C8").NumberFormat = $ws.Range("C9").NumberFormat
